$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 0.9487032799389779
$ws.Range("C3").Value = 0.943078634153075
$ws.Range("D3").Value = 0.9487032799389779
$ws.Range("E3").Value = 0.9383652640608233

$ws.Range("B4").Value = 0.0482280533236697
$ws.Range("C4").Value = 0.05060722202633464
$ws.Range("D4").Value = 0.0482280533236697
$ws.Range("E4").Value = 0.05126769446794536

$ws.Range("B5").Value = 0.7864225781845919
$ws.Range("C5").Value = 0.7872906718385223
$ws.Range("D5").Value = 0.7864225781845919
$ws.Range("E5").Value = 0.7831637240824391

$ws.Range("B6").Value = 0.9429824561403508
$ws.Range("C6").Value = 0.9360374286231316
$ws.Range("D6").Value = 0.9429824561403508
$ws.Range("E6").Value = 0.9300699916820448

$ws.Range("B7").Value = 0.9603356216628528
$ws.Range("C7").Value = 0.9576300588995028
$ws.Range("D7").Value = 0.9603356216628528
$ws.Range("E7").Value = 0.9496149427735556

$ws.Range("B8").Value = 0.9802631578947368
$ws.Range("C8").Value = 0.9776320668894856
$ws.Range("D8").Value = 0.9802631578947368
$ws.Range("E8").Value = 0.9727968902495393

$ws.Range("B9").Value = 0.9956140350877193
$ws.Range("C9").Value = 0.9956332717759311
$ws.Range("D9").Value = 0.9956140350877193
$ws.Range("E9").Value = 0.9934258723732408
